# Edit test_expense_data.xlsx to match the obj_tables-style "ObjTables" header
# convention: insert a marker row above the existing header row, prefix each
# header label with "!", and rename the sheet to "!!Transactions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Insert a new blank row above the current header row (row 1), pushing all
# existing rows (header + data) down by one.
$ws.Rows.Item(1).Insert()

# New row 1: the ObjTables table-declaration marker, in column A only, with
# plain/default formatting (no inherited style from the row below).
$ws.Range("A1").Value = "!!ObjTables type='Data' class='Transaction' tableFormat='row'"
$ws.Range("A1").Style = "Normal"

# The former header row is now row 2: prefix each header label with "!"
# (columns C "Extra column" and G "Another extra" are left as-is).
$ws.Range("A2").Value = "!Date"
$ws.Range("B2").Value = "!Payee"
$ws.Range("D2").Value = "!amount"
$ws.Range("E2").Value = "!Tax category"
$ws.Range("F2").Value = "!spending_category"

# Rename the worksheet to mark it as an ObjTables sheet too.
$ws.Name = "!!Transactions"

# Restore view state: zoom 120% and selection on F3 (was D7 at 130%).
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
[void]$ws.Range("F3").Select()
